$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# covering 2021-05-28 through 2021-06-28 (aggiornamento fino a 28/06 incluso)
$data = @(
    @(44344, 0, 1, 9.930486593843098),
    @(44345, 1, 2, 19.8609731876862),
    @(44346, 1, 3, 29.7914597815293),
    @(44347, 1, 4, 39.72194637537239),
    @(44348, 0, 4, 39.72194637537239),
    @(44349, 0, 3, 29.7914597815293),
    @(44350, 0, 3, 29.7914597815293),
    @(44351, 0, 3, 29.7914597815293),
    @(44352, 0, 2, 19.8609731876862),
    @(44353, 0, 1, 9.930486593843098),
    @(44354, 0, 0, 0),
    @(44355, 0, 0, 0),
    @(44356, 0, 0, 0),
    @(44357, 0, 0, 0),
    @(44358, 0, 0, 0),
    @(44359, 0, 0, 0),
    @(44360, 0, 0, 0),
    @(44361, 0, 0, 0),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 0, 0, 0),
    @(44365, 0, 0, 0),
    @(44366, 0, 0, 0),
    @(44367, 0, 0, 0),
    @(44368, 0, 0, 0),
    @(44369, 0, 0, 0),
    @(44370, 0, 0, 0),
    @(44371, 0, 0, 0),
    @(44372, 0, 0, 0),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$startRow = 270
$endRow = $startRow + $data.Count - 1

# Copy the formatting of the last existing row (269) down across all the new
# rows at once, so the date column keeps its date style/border/etc. and the
# other columns keep the default (unstyled) look.
$ws.Range("A269:D269").Copy() | Out-Null
$ws.Range("A$startRow`:D$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
